$wb = $excel.ActiveWorkbook

# --- Training Dashboard: row 3 "PERIOD TO EXPIRE" / "LAST UPDATE" refresh ---
$wsTraining = $wb.Worksheets.Item("Training Dashboard")
$wsTraining.Range("H3").Value = -343

# Keep the refreshed date a literal text value (matches how the sheet already
# stores its other date-looking cells) instead of letting it auto-convert to
# a date serial.
$dateCell = $wsTraining.Range("I3")
$dateCell.NumberFormat = "@"
$dateCell.Value = "16-Sep-2025"

# --- Exam Dashboard: comments column narrower + remarks now say the date is fine ---
$wsExam = $wb.Worksheets.Item("Exam Dashboard")
$wsExam.Columns.Item(5).ColumnWidth = 15 - 0.83

$wsExam.Range("E3").Value = "date is valid"
$wsExam.Range("E4").Value = "date is valid"

# --- Header styling: bold header text becomes white on the dark-blue fill ---
$wsTraining.Range("A2:K2").Font.Color = 16777215
$wsExam.Range("A2:G2").Font.Color = 16777215
